# Rename all 30 worksheets in the workbook, keeping their order, sheetId and
# r:id the same, but updating the visible sheet name (name-only change -
# the content of each sheet is left untouched).
#
# Mapping is positional: the Nth sheet (in the existing sheet order) takes
# on the Nth new name below.

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ23603103",
    "summ23759568",
    "summ23954879",
    "summ24151174",
    "summ24363321",
    "summ24573365",
    "summ24792742",
    "summ24995431",
    "summ25198882",
    "summ25406729",
    "summ25611764",
    "summ25807079",
    "summ26040615",
    "summ26256372",
    "summ26450246",
    "summ26666152",
    "summ26877313",
    "summ27091504",
    "summ27288308",
    "summ27563708",
    "summ27778626",
    "summ28058547",
    "summ28257650",
    "summ28456288",
    "summ28683445",
    "summ28914000",
    "summ29125567",
    "summ29349923",
    "summ29563949",
    "summ29773984"
)

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = $newNames[$i]
}
